$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.011.92"
$ws.Range("E2").Value = "  +2.33%  "

$ws.Range("D3").Value = "2.048.25"
$ws.Range("E3").Value = "  +1.44%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "228.73"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("D7").Value = "60.52"
$ws.Range("E7").Value = "  +8.16%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +2.19%  "

$ws.Range("E10").Value = "  +2.93%  "

$ws.Range("E11").Value = "  +1.81%  "

$ws.Range("D12").Value = "14.71"
$ws.Range("E12").Value = "  +3.18%  "

$ws.Range("D13").Value = "2.352.06"
$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("D14").Value = "20.96"
$ws.Range("E14").Value = "  +5.00%  "

$ws.Range("D15").Value = "5.31"
$ws.Range("E15").Value = "  +2.73%  "

$ws.Range("E16").Value = "  +2.51%  "

$ws.Range("D17").Value = "2.050.80"
$ws.Range("E17").Value = "  +1.74%  "

$ws.Range("D18").Value = "37.985.07"
$ws.Range("E18").Value = "  +2.56%  "

$ws.Range("D19").Value = "6.29"
$ws.Range("E19").Value = "  +2.54%  "

$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("D21").Value = "0.0@@SUB3@@0829"
$ws.Range("D21").Replace("@@SUB3@@", [char]0x2083) | Out-Null
$ws.Range("E21").Value = "  +1.74%  "

$ws.Range("D22").Value = "225.52"
$ws.Range("E22").Value = "  +1.23%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("E25").Value = "  +1.48%  "

$ws.Range("D26").Value = "165.63"
$ws.Range("E26").Value = "  +1.42%  "

$ws.Range("D27").Value = "9.17"
$ws.Range("E27").Value = "  +1.93%  "

$ws.Range("D28").Value = "0.132"
$ws.Range("E28").Value = "  +3.70%  "

$ws.Range("D29").Value = "18.99"
$ws.Range("E29").Value = "  +1.58%  "

$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("D31").Value = "0.119"
$ws.Range("E31").Value = "  +1.90%  "

$ws.Range("E32").Value = "  +1.60%  "

$ws.Range("D33").Value = "4.55"
$ws.Range("E33").Value = "  +2.21%  "

$ws.Range("E34").Value = "  +8.51%  "

$ws.Range("D36").Value = "6.24"
$ws.Range("E36").Value = "  +14.22%  "

$ws.Range("D37").Value = "2.29"
$ws.Range("E37").Value = "  -2.46%  "

$ws.Range("E38").Value = "  +2.59%  "

$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").Value = "1.518.69"
$ws.Range("E40").Value = "  +3.34%  "

$ws.Range("D41").Value = "97.25"
$ws.Range("E41").Value = "  +3.57%  "

$ws.Range("D42").Value = "16.86"
$ws.Range("E42").Value = "  +3.53%  "

$ws.Range("E43").Value = "  +1.23%  "

$ws.Range("E44").Value = "  +2.72%  "

$ws.Range("D45").Value = "0.0921"
$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("D46").Value = "1.12"
$ws.Range("E46").Value = "  +1.76%  "

$ws.Range("D47").Value = "4.02"
$ws.Range("E47").Value = "  -6.75%  "

$ws.Range("E48").Value = "  +0.62%  "

$ws.Range("E49").Value = "  +1.49%  "

$ws.Range("D50").Value = "7.01"
$ws.Range("E50").Value = "  -0.42%  "

$ws.Range("D51").Value = "2.241.03"
